$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing week-number lists to include the new week "43" ---
# (these are the only real content changes among the pre-existing rows;
#  the remaining shared-string index churn in the diff is just a side
#  effect of the string table shrinking/growing and needs no extra edits)
$ws.Range("D3").Value  = "15, 19, 21, 23, 35, 43"
$ws.Range("D7").Value  = "35, 39, 43"
$ws.Range("D8").Value  = "21, 23, 25, 27, 35, 39, 43"
$ws.Range("D11").Value = "35, 39, 43"
$ws.Range("D12").Value = "33, 35, 41, 43"
$ws.Range("D13").Value = "21, 23, 25, 27, 29, 31, 33, 35, 37, 39, 43"

# --- Add the new "Stimulus" module row (row 19) ---
# Match formatting of the rest of the table for the new row first (it mirrors
# a normal, non-bordered data row like row 18), then fill in the values.
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)
$ws.Range("A19").RowHeight = $ws.Range("A18").RowHeight

$ws.Range("A19").Value = "Stimulus"
$ws.Range("B19").Value = "Did you receive stimulus?; How did you spend stimulus? "
$ws.Range("C19").Value = "https://r3questionbank.netlify.app/stimulusmodule"
$ws.Range("D19").Value = 43

# Hyperlink the new module's URL cell, then restore its non-hyperlink-blue
# look to match the rest of column C (which uses plain formatting, not the
# auto Hyperlink style).
$ws.Hyperlinks.Add($ws.Range("C19"), "https://r3questionbank.netlify.app/stimulusmodule")
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122)

# --- Cosmetic: selection, matching the saved view ---
$ws.Range("D19").Select() | Out-Null
